$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose new values would otherwise be
# auto-converted to numbers (losing trailing zeros / becoming scientific notation),
# matching the original inline-string text storage of the Price column.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "64.193.18"
$ws.Range("E2").Value = "  +1.75%  "
$ws.Range("D3").Value = "2.675.51"
$ws.Range("E3").Value = "  +3.08%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "594.63"
$ws.Range("E5").Value = "  +1.98%  "
$ws.Range("D6").Value = "147.43"
$ws.Range("E6").Value = "  -0.55%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("E8").Value = "  -0.91%  "
$ws.Range("E9").Value = "  +0.02%  "
$ws.Range("D10").Value = "5.65"
$ws.Range("E10").Value = "  -0.47%  "
$ws.Range("E11").Value = "  -0.15%  "
$ws.Range("E12").Value = "  +0.55%  "
$ws.Range("D13").Value = "27.78"
$ws.Range("E13").Value = "  +1.80%  "
$ws.Range("D14").Value = "3.154.84"
$ws.Range("E14").Value = "  +3.09%  "
$ws.Range("D15").Value = "64.176.84"
$ws.Range("E15").Value = "  +1.91%  "
$ws.Range("D16").Value = "0.0000147"
$ws.Range("E16").Value = "  -0.11%  "
$ws.Range("D17").Value = "2.652.77"
$ws.Range("E17").Value = "  +1.96%  "
$ws.Range("D18").Value = "11.41"
$ws.Range("E18").Value = "  +0.28%  "
$ws.Range("D19").Value = "345.17"
$ws.Range("E19").Value = "  +0.30%  "
$ws.Range("E20").Value = "  -0.39%  "
$ws.Range("D21").Value = "6.87"
$ws.Range("E21").Value = "  +1.15%  "
$ws.Range("E22").Value = "  +0.27%  "
$ws.Range("D23").Value = "68.59"
$ws.Range("E23").Value = "  +1.99%  "
$ws.Range("E24").Value = "  +10.71%  "
$ws.Range("D25").Value = "1.66"
$ws.Range("E25").Value = "  +3.70%  "
$ws.Range("E26").Value = "  -1.61%  "
$ws.Range("B27").Value = "Bittensor"
$ws.Range("C27").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D27").Value = "547.18"
$ws.Range("E27").Value = "  +15.95%  "
$ws.Range("B28").Value = "InternetComputer(DFINITY)"
$ws.Range("C28").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D28").Value = "8.54"
$ws.Range("E28").Value = "  +1.62%  "
$ws.Range("B29").Value = "Aptos"
$ws.Range("C29").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D29").Value = "7.98"
$ws.Range("E29").Value = "  +1.16%  "
$ws.Range("B30").Value = "Binance-PegBSC-USD"
$ws.Range("C30").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D30").Value = "1.00"
$ws.Range("E30").Value = "  +0.27%  "
$ws.Range("D31").Value = "2.01"
$ws.Range("E31").Value = "  +3.21%  "
$ws.Range("D32").Value = "1.79"
$ws.Range("E32").Value = "  +11.05%  "
$ws.Range("D33").Value = "0.0₃0823"
$ws.Range("E33").Value = "  -0.25%  "
$ws.Range("D34").Value = "175.18"
$ws.Range("E34").Value = "  -0.88%  "
$ws.Range("E35").Value = "  +0.10%  "
$ws.Range("E36").Value = "  -0.81%  "
$ws.Range("D37").Value = "19.30"
$ws.Range("E37").Value = "  +0.44%  "
$ws.Range("D38").Value = "4.72"
$ws.Range("E38").Value = "  +2.77%  "
$ws.Range("B39").Value = "Stacks"
$ws.Range("C39").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D39").Value = "1.76"
$ws.Range("E39").Value = "  +2.99%  "
$ws.Range("B40").Value = "Aave"
$ws.Range("C40").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D40").Value = "172.94"
$ws.Range("E40").Value = "  +8.65%  "
$ws.Range("E41").Value = "  +0.03%  "
$ws.Range("D42").Value = "40.57"
$ws.Range("E42").Value = "  +2.70%  "
$ws.Range("D43").Value = "3.78"
$ws.Range("E43").Value = "  -0.57%  "
$ws.Range("D44").Value = "21.59"
$ws.Range("E44").Value = "  +1.81%  "
$ws.Range("D45").Value = "0.635"
$ws.Range("E45").Value = "  -0.28%  "
$ws.Range("D46").Value = "0.0549"
$ws.Range("E46").Value = "  -0.12%  "
$ws.Range("E47").Value = "  +0.86%  "
$ws.Range("D48").Value = "0.0965"
$ws.Range("E48").Value = "  -0.96%  "
$ws.Range("E49").Value = "  +1.23%  "
$ws.Range("D50").Value = "1.77"
$ws.Range("E50").Value = "  +1.94%  "
$ws.Range("D51").Value = "11.34"
$ws.Range("E51").Value = "  -1.01%  "
